$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Experiment Number values in rows 5-8 ---
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6

# --- Add result text in K7/K8 (reusing existing shared string) ---
$ws.Range("K7").Value = "micro corr and macro corr condition did not meet"
$ws.Range("K8").Value = "micro corr and macro corr condition did not meet"

# --- Fill in column K (results) first for rows 9-10 so new shared strings are
#     registered in the same order as the target workbook ---
$ws.Range("K9").Value = "micro corelation min is greater than macro corelation max"
$ws.Range("K10").Value = "micro corelation min is greater than macro corelation max (Current best output)"

# --- Fill the remaining data columns (A-I) for the four new experiment rows ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 13
$ws.Range("E9").Value = 0.5
$ws.Range("F9").Value = 0.2
$ws.Range("G9").Value = -1
$ws.Range("H9").Value = "32x32"
$ws.Range("I9").Value = "64x64"

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = 13
$ws.Range("E10").Value = 0.5
$ws.Range("F10").Value = 0.3
$ws.Range("G10").Value = -1
$ws.Range("H10").Value = "32x32"
$ws.Range("I10").Value = "64x64"

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 13
$ws.Range("E11").Value = 0.5
$ws.Range("F11").Value = 0.35
$ws.Range("G11").Value = -1
$ws.Range("H11").Value = "32x32"
$ws.Range("I11").Value = "64x64"

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = $false
$ws.Range("D12").Value = 13
$ws.Range("E12").Value = 0.5
$ws.Range("F12").Value = 0.4
$ws.Range("G12").Value = -1
$ws.Range("H12").Value = "32x32"
$ws.Range("I12").Value = "64x64"

# --- Add the URL hyperlinks for the new rows (J9-J12), registering the
#     remaining new shared strings in order ---
$ws.Hyperlinks.Add($ws.Range("J9"), "https://github.com/GurunagSai/neocortexapi-classification/blob/GurunagSai/ExperimentReport/DataSet-3/output-07.png", "", "", "https://github.com/GurunagSai/neocortexapi-classification/blob/GurunagSai/ExperimentReport/DataSet-3/output-07.png")
$ws.Range("J9").WrapText = $true

$ws.Hyperlinks.Add($ws.Range("J10"), "https://github.com/GurunagSai/neocortexapi-classification/blob/GurunagSai/ExperimentReport/DataSet-3/output-08.png", "", "", "https://github.com/GurunagSai/neocortexapi-classification/blob/GurunagSai/ExperimentReport/DataSet-3/output-08.png")
$ws.Range("J10").WrapText = $true

$ws.Hyperlinks.Add($ws.Range("J11"), "https://github.com/GurunagSai/neocortexapi-classification/blob/GurunagSai/ExperimentReport/DataSet-3/output-09.png", "", "", "https://github.com/GurunagSai/neocortexapi-classification/blob/GurunagSai/ExperimentReport/DataSet-3/output-09.png")
$ws.Range("J11").WrapText = $true

$ws.Hyperlinks.Add($ws.Range("J12"), "https://github.com/GurunagSai/neocortexapi-classification/blob/GurunagSai/ExperimentReport/DataSet-3/output-09.png", "", "", "https://github.com/GurunagSai/neocortexapi-classification/blob/GurunagSai/ExperimentReport/DataSet-3/output-09.png")
$ws.Range("J12").WrapText = $true

# --- K11/K12 reuse the existing "did not meet" text ---
$ws.Range("K11").Value = "micro corr and macro corr condition did not meet"
$ws.Range("K12").Value = "micro corr and macro corr condition did not meet"

# --- Update the view: select H17 (clears the old C1 top-left/ K6 selection) ---
$ws.Range("H17").Select()
